$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Прибор поверки" value for row 2 (replaces the placeholder blank spaces)
$ws.Range("I2").Value = "Установка для поверки каналов измерения давления и частоты пульса УПКД-3"

# New column L ("range") with its header and first four data values
$ws.Range("L1").Value = "range"
$ws.Range("L2").Value = 1800
$ws.Range("L3").Value = 200
$ws.Range("L4").Value = 60
$ws.Range("L5").Value = 40

# Update view: scroll/selection moved to the newly added column
$ws.Range("L3").Select()
